$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.25668576553308
$ws.Range("C2").Value = 8.294962618303064
$ws.Range("D2").Value = 6.016425161195423
$ws.Range("E2").Value = 10.63261176225111
$ws.Range("G2").Value = 55.14396153613016
$ws.Range("H2").Value = 20.06463767247139
$ws.Range("L2").Value = 10.15134418404072
$ws.Range("B3").Value = 19.84175710689192
$ws.Range("C3").Value = 7.897022757154873
$ws.Range("D3").Value = 5.904646904057813
$ws.Range("E3").Value = 10.65024970330101
$ws.Range("G3").Value = 54.30185226203805
$ws.Range("H3").Value = 19.98838200967699
$ws.Range("L3").Value = 10.13800711783105
$ws.Range("B4").Value = 19.58960932346353
$ws.Range("C4").Value = 7.640293370363234
$ws.Range("D4").Value = 5.836955220237453
$ws.Range("E4").Value = 10.66216224908986
$ws.Range("G4").Value = 53.7948506004724
$ws.Range("H4").Value = 19.9459583214288
$ws.Range("L4").Value = 10.13202632496074
$ws.Range("B5").Value = 19.48767993894712
$ws.Range("C5").Value = 7.532614652266297
$ws.Range("D5").Value = 5.809646447545986
$ws.Range("E5").Value = 10.66728877307331
$ws.Range("G5").Value = 53.59100702358845
$ws.Range("H5").Value = 19.92978142798462
$ws.Range("L5").Value = 10.13014561440271
$ws.Range("B6").Value = 19.47080915925645
$ws.Range("C6").Value = 7.51455177742906
$ws.Range("D6").Value = 5.805129678281349
$ws.Range("E6").Value = 10.66815645435867
$ws.Range("G6").Value = 53.55733260610347
$ws.Range("H6").Value = 19.92716250155535
$ws.Range("L6").Value = 10.12986696400514
$ws.Range("B7").Value = 19.58823111876956
$ws.Range("C7").Value = 7.638853481806922
$ws.Range("D7").Value = 5.836585754268219
$ws.Range("E7").Value = 10.6622302858734
$ws.Range("G7").Value = 53.79208999933023
$ws.Range("H7").Value = 19.94573564963491
$ws.Range("L7").Value = 10.13199870642287
$ws.Range("B8").Value = 20.1131724747924
$ws.Range("C8").Value = 8.160355915503265
$ws.Range("D8").Value = 5.977712586586239
$ws.Range("E8").Value = 10.63846846926343
$ws.Range("G8").Value = 54.85165770685717
$ws.Range("H8").Value = 20.03743559285706
$ws.Range("L8").Value = 10.14628810517802
$ws.Range("B9").Value = 21.15627642076233
$ws.Range("C9").Value = 9.082701905637318
$ws.Range("D9").Value = 6.260188039079089
$ws.Range("E9").Value = 10.60047249549526
$ws.Range("G9").Value = 56.99869394584235
$ws.Range("H9").Value = 20.25185831349845
$ws.Range("L9").Value = 10.1917604400492
$ws.Range("B10").Value = 21.92177709249442
$ws.Range("C10").Value = 9.697132702895352
$ws.Range("D10").Value = 6.468922681686521
$ws.Range("E10").Value = 10.57781349160549
$ws.Range("G10").Value = 58.60400129078809
$ws.Range("H10").Value = 20.43000850298402
$ws.Range("L10").Value = 10.23568992052234
$ws.Range("B11").Value = 22.2679978087052
$ws.Range("C11").Value = 9.962631666667662
$ws.Range("D11").Value = 6.56367581936895
$ws.Range("E11").Value = 10.56864958763318
$ws.Range("G11").Value = 59.33732105447195
$ws.Range("H11").Value = 20.51539998857455
$ws.Range("L11").Value = 10.25792575531729
$ws.Range("B12").Value = 22.39867053296141
$ws.Range("C12").Value = 10.06113704074334
$ws.Range("D12").Value = 6.599490454199526
$ws.Range("E12").Value = 10.56534416174989
$ws.Range("G12").Value = 59.61520139521167
$ws.Range("H12").Value = 20.54834798759135
$ws.Range("L12").Value = 10.26666622311337
$ws.Range("B13").Value = 22.3705493194045
$ws.Range("C13").Value = 10.040012831098
$ws.Range("D13").Value = 6.591780658115949
$ws.Range("E13").Value = 10.5660487130038
$ws.Range("G13").Value = 59.55535046602974
$ws.Range("H13").Value = 20.5412250497516
$ws.Range("L13").Value = 10.26476962585001
$ws.Range("B14").Value = 22.27875781991952
$ws.Range("C14").Value = 9.970776605874995
$ws.Range("D14").Value = 6.566623848044416
$ws.Range("E14").Value = 10.56837434497753
$ws.Range("G14").Value = 59.36018008694853
$ws.Range("H14").Value = 20.51809845412576
$ws.Range("L14").Value = 10.25863844284038
$ws.Range("B15").Value = 22.2224722658466
$ws.Range("C15").Value = 9.928102137885743
$ws.Range("D15").Value = 6.551204858596214
$ws.Range("E15").Value = 10.56982032439535
$ws.Range("G15").Value = 59.24064972521022
$ws.Range("H15").Value = 20.50401204056157
$ws.Range("L15").Value = 10.25492450710225
$ws.Range("B16").Value = 21.89909777154073
$ws.Range("C16").Value = 9.67949810929891
$ws.Range("D16").Value = 6.462723056915842
$ws.Range("E16").Value = 10.57843541345963
$ws.Range("G16").Value = 58.55611699849309
$ws.Range("H16").Value = 20.42451438163734
$ws.Range("L16").Value = 10.23428177481661
$ws.Range("B17").Value = 21.70009786389802
$ws.Range("C17").Value = 9.523384682573115
$ws.Range("D17").Value = 6.408363655295359
$ws.Range("E17").Value = 10.58401363735176
$ws.Range("G17").Value = 58.13678056328501
$ws.Range("H17").Value = 20.37685012085921
$ws.Range("L17").Value = 10.22219250020277
$ws.Range("B18").Value = 21.58545739460933
$ws.Range("C18").Value = 9.432274609103814
$ws.Range("D18").Value = 6.377080753013968
$ws.Range("E18").Value = 10.58732972736269
$ws.Range("G18").Value = 57.89589184514629
$ws.Range("H18").Value = 20.34984523095169
$ws.Range("L18").Value = 10.21545128829836
$ws.Range("B19").Value = 21.5466155574767
$ws.Range("C19").Value = 9.401200630016172
$ws.Range("D19").Value = 6.366487178797128
$ws.Range("E19").Value = 10.58847097796993
$ws.Range("G19").Value = 57.81439078300852
$ws.Range("H19").Value = 20.34077270390864
$ws.Range("L19").Value = 10.21320538209195
$ws.Range("B20").Value = 21.72130154887223
$ws.Range("C20").Value = 9.540139757946655
$ws.Range("D20").Value = 6.414152322566021
$ws.Range("E20").Value = 10.58340868326862
$ws.Range("G20").Value = 58.18139017604523
$ws.Range("H20").Value = 20.38188168853782
$ws.Range("L20").Value = 10.22345748514905
$ws.Range("B21").Value = 22.3057320857661
$ws.Range("C21").Value = 9.991168275263588
$ws.Range("D21").Value = 6.57401510221282
$ws.Range("E21").Value = 10.56768677721569
$ws.Range("G21").Value = 59.41750322809403
$ws.Range("H21").Value = 20.52487479665753
$ws.Range("L21").Value = 10.26043065939395
$ws.Range("B22").Value = 22.68510309035037
$ws.Range("C22").Value = 10.27408404224543
$ws.Range("D22").Value = 6.678092015643617
$ws.Range("E22").Value = 10.55837198849371
$ws.Range("G22").Value = 60.22635881459899
$ws.Range("H22").Value = 20.62188977916603
$ws.Range("L22").Value = 10.28645942145894
$ws.Range("B23").Value = 22.48290805575629
$ws.Range("C23").Value = 10.12417662384642
$ws.Range("D23").Value = 6.622592943523515
$ws.Range("E23").Value = 10.56325550320209
$ws.Range("G23").Value = 59.7946487193881
$ws.Range("H23").Value = 20.56978999092454
$ws.Range("L23").Value = 10.27239805765487
$ws.Range("B24").Value = 21.71171607686427
$ws.Range("C24").Value = 9.532569017955648
$ws.Range("D24").Value = 6.411535359666
$ws.Range("E24").Value = 10.58368184322423
$ws.Range("G24").Value = 58.16122154861265
$ws.Range("H24").Value = 20.37960567632684
$ws.Range("L24").Value = 10.22288493396963
$ws.Range("B25").Value = 20.87359514000309
$ws.Range("C25").Value = 8.844164125064346
$ws.Range("D25").Value = 6.183397690328204
$ws.Range("E25").Value = 10.60982923827612
$ws.Range("G25").Value = 56.41192963636967
$ws.Range("H25").Value = 20.19019118727472
$ws.Range("L25").Value = 10.17760124361842
